# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF with labels, matching the style of
# the existing header row (bold font, border, centered/top aligned).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every data row (2-37): this team went 92-70-0.
$wins = 92
$losses = 70
$ties = 0

for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
